{"js": "// Update invasive species classification table: replace the rounded\n// 2-decimal statistic values (Global Mean, Lower/Higher Credible Interval)\n// with their full 4-decimal precision values.\n//\n// The replacements are applied in table (row-major, then column) order,\n// which is also the order the values appear in the document body, so a\n// simple \"next value in the queue\" lookup keyed by the old text correctly\n// disambiguates duplicate old values (e.g. \"0.00\" / \"-0.05\" / \"-0.18\"\n// appear more than once).\n\nconst replacements = [\n  [\"1.09\", \"1.0920\"],\n  [\"-2.75\", \"-2.7486\"],\n  [\"4.37\", \"4.3651\"],\n  [\"0.73\", \"0.7316\"],\n  [\"-3.65\", \"-3.6508\"],\n  [\"4.56\", \"4.5637\"],\n  [\"-0.05\", \"-0.0453\"],\n  [\"-0.18\", \"-0.1822\"],\n  [\"0.09\", \"0.0868\"],\n  [\"-0.07\", \"-0.0692\"],\n  [\"-0.29\", \"-0.2945\"],\n  [\"0.12\", \"0.1189\"],\n  [\"-3.90\", \"-3.9049\"],\n  [\"-8.94\", \"-8.9421\"],\n  [\"1.79\", \"1.7926\"],\n  [\"-13.07\", \"-13.0663\"],\n  [\"-17.53\", \"-17.5297\"],\n  [\"-8.35\", \"-8.3527\"],\n  [\"0.66\", \"0.6607\"],\n  [\"-14.88\", \"-14.8789\"],\n  [\"16.70\", \"16.6952\"],\n  [\"-2.27\", \"-2.2745\"],\n  [\"-22.75\", \"-22.7494\"],\n  [\"19.15\", \"19.1467\"],\n  [\"-0.73\", \"-0.7266\"],\n  [\"-1.89\", \"-1.8935\"],\n  [\"0.33\", \"0.3324\"],\n  [\"-1.88\", \"-1.8755\"],\n  [\"-2.86\", \"-2.8622\"],\n  [\"-0.88\", \"-0.8787\"],\n  [\"1.37\", \"1.3723\"],\n  [\"-1.16\", \"-1.1642\"],\n  [\"3.48\", \"3.4838\"],\n  [\"1.41\", \"1.4058\"],\n  [\"-1.86\", \"-1.8569\"],\n  [\"3.85\", \"3.8460\"],\n  [\"0.08\", \"0.0783\"],\n  [\"0.00\", \"-0.0017\"],\n  [\"0.16\", \"0.1616\"],\n  [\"0.00\", \"0.0015\"],\n  [\"-0.05\", \"-0.0468\"],\n  [\"0.05\", \"0.0518\"],\n  [\"0.03\", \"0.0250\"],\n  [\"-0.16\", \"-0.1642\"],\n  [\"0.21\", \"0.2052\"],\n  [\"0.06\", \"0.0583\"],\n  [\"-0.18\", \"-0.1794\"],\n  [\"0.30\", \"0.2952\"],\n];\n\n// Index replacements by the old text, preserving the order each distinct\n// old value appears so repeated values are consumed in sequence.\nconst queues = new Map();\nfor (const [oldText, newText] of replacements) {\n  if (!queues.has(oldText)) queues.set(oldText, []);\n  queues.get(oldText).push(newText);\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst values = table.values;\nconst rowCount = values.length;\n\n// The three rightmost-but-one numeric columns (\"Global Mean\", \"Lower\n// Credible Interval\", \"Higher Credible Interval\") are columns 3, 4, 5\n// (0-indexed) of the 7-column table; the last column is the units text.\nconst dataColumns = [3, 4, 5];\n\nfor (let r = 0; r < rowCount; r++) {\n  for (const c of dataColumns) {\n    const cellText = values[r][c];\n    const queue = queues.get(cellText);\n    if (queue && queue.length > 0) {\n      const newText = queue.shift();\n      const cell = table.getCell(r, c);\n      // Use the cell's Range (scoped to the existing paragraph/run) rather\n      // than cell.body, so only the text node is rewritten and the\n      // surrounding run/paragraph formatting (rPr/pPr) is preserved.\n      const range = cell.getRange();\n      range.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update invasive species classification table: replace the rounded\n# 2-decimal statistic values (Global Mean, Lower/Higher Credible Interval)\n# with their full 4-decimal precision values.\n#\n# The replacements are applied in table (row-major, then column) order,\n# which is also the order the values appear in the document, so a simple\n# \"next value in the queue\" lookup keyed by the old text correctly\n# disambiguates duplicate old values (e.g. \"0.00\" / \"-0.05\" / \"-0.18\"\n# appear more than once).\n\n$replacements = @(\n  @(\"1.09\", \"1.0920\"),\n  @(\"-2.75\", \"-2.7486\"),\n  @(\"4.37\", \"4.3651\"),\n  @(\"0.73\", \"0.7316\"),\n  @(\"-3.65\", \"-3.6508\"),\n  @(\"4.56\", \"4.5637\"),\n  @(\"-0.05\", \"-0.0453\"),\n  @(\"-0.18\", \"-0.1822\"),\n  @(\"0.09\", \"0.0868\"),\n  @(\"-0.07\", \"-0.0692\"),\n  @(\"-0.29\", \"-0.2945\"),\n  @(\"0.12\", \"0.1189\"),\n  @(\"-3.90\", \"-3.9049\"),\n  @(\"-8.94\", \"-8.9421\"),\n  @(\"1.79\", \"1.7926\"),\n  @(\"-13.07\", \"-13.0663\"),\n  @(\"-17.53\", \"-17.5297\"),\n  @(\"-8.35\", \"-8.3527\"),\n  @(\"0.66\", \"0.6607\"),\n  @(\"-14.88\", \"-14.8789\"),\n  @(\"16.70\", \"16.6952\"),\n  @(\"-2.27\", \"-2.2745\"),\n  @(\"-22.75\", \"-22.7494\"),\n  @(\"19.15\", \"19.1467\"),\n  @(\"-0.73\", \"-0.7266\"),\n  @(\"-1.89\", \"-1.8935\"),\n  @(\"0.33\", \"0.3324\"),\n  @(\"-1.88\", \"-1.8755\"),\n  @(\"-2.86\", \"-2.8622\"),\n  @(\"-0.88\", \"-0.8787\"),\n  @(\"1.37\", \"1.3723\"),\n  @(\"-1.16\", \"-1.1642\"),\n  @(\"3.48\", \"3.4838\"),\n  @(\"1.41\", \"1.4058\"),\n  @(\"-1.86\", \"-1.8569\"),\n  @(\"3.85\", \"3.8460\"),\n  @(\"0.08\", \"0.0783\"),\n  @(\"0.00\", \"-0.0017\"),\n  @(\"0.16\", \"0.1616\"),\n  @(\"0.00\", \"0.0015\"),\n  @(\"-0.05\", \"-0.0468\"),\n  @(\"0.05\", \"0.0518\"),\n  @(\"0.03\", \"0.0250\"),\n  @(\"-0.16\", \"-0.1642\"),\n  @(\"0.21\", \"0.2052\"),\n  @(\"0.06\", \"0.0583\"),\n  @(\"-0.18\", \"-0.1794\"),\n  @(\"0.30\", \"0.2952\")\n)\n\n# Index replacements by the old text, preserving the order each distinct\n# old value appears so repeated values are consumed in sequence. Each\n# queue is a plain array plus a \"next index\" cursor (avoids relying on\n# .NET generic collection support in the PowerShell host).\n$queues = @{}\n$cursors = @{}\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  if (-not $queues.ContainsKey($oldText)) {\n    $queues[$oldText] = @()\n    $cursors[$oldText] = 0\n  }\n  $queues[$oldText] += $newText\n}\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# The three rightmost-but-one numeric columns (\"Global Mean\", \"Lower\n# Credible Interval\", \"Higher Credible Interval\") are columns 4, 5, 6\n# (1-indexed) of the 7-column table; the last column is the units text.\n$dataColumns = @(4, 5, 6)\n\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n  foreach ($c in $dataColumns) {\n    $cell = $table.Cell($r, $c)\n    $range = $cell.Range\n    # Cell range text includes trailing cell-mark characters (CR + BEL);\n    # strip them to get the plain cell value for lookup.\n    $cellText = $range.Text.TrimEnd([char]13, [char]7)\n    if ($queues.ContainsKey($cellText) -and $cursors[$cellText] -lt $queues[$cellText].Count) {\n      $newText = $queues[$cellText][$cursors[$cellText]]\n      $cursors[$cellText] = $cursors[$cellText] + 1\n      # Replace only the literal text portion of the range (excluding the\n      # trailing cell-end mark) so run/paragraph formatting is preserved.\n      $textRange = $d.Range($range.Start, $range.Start + $cellText.Length)\n      $textRange.Text = $newText\n    }\n  }\n}\n"}
